$wb = $excel.ActiveWorkbook

# --- Unhide the "Sheet2" lookup sheet ---
$lookup = $wb.Worksheets.Item("Sheet2")
$lookup.Visible = -1

# --- Update the SOURCE_TYPE lookup list on "Sheet2" (column I) ---
# Was: I1=Ex Situ Excavation, I2=In Situ Excavation, I3=Other
# Now: I1=Stockpile, I2=Excavation, I3=Drill Spoils, I4=Other (list grew by one entry)
$lookup.Range("I1").Value = "Stockpile"
$lookup.Range("I2").Value = "Excavation"
$lookup.Range("I3").Value = "Drill Spoils"
$lookup.Range("I4").Value = "Other"

# --- Update the data validation source range on the "Create Batch" sheet (H column) ---
$batch = $wb.Worksheets.Item("Create Batch")
$batch.Range("H2:H1048576").Validation.Delete()
$batch.Range("H2:H1048576").Validation.Add(3, 1, 1, "=Sheet2!`$I`$1:`$I`$4")
$batch.Range("H2:H1048576").Validation.InCellDropdown = $true

# --- Update the sample/test data row on "Create Batch" ---
$batch.Range("B2").Value = "10/25/2020"
$batch.Range("C2").Value = "TESTTEST10252020"
$batch.Range("H2").Value = "Stockpile"

# --- Restore view/selection state ---
$lookup.Range("C22").Select()
$batch.Activate()
$batch.Range("C3").Select()
